$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "36÷8=4, 4" "35÷5=7, 0"
Replace-Text "67÷2=33, 1" "16÷6=2, 4"
Replace-Text "11÷4=2, 3" "20÷4=5, 0"
Replace-Text "75÷6=12, 3" "44÷2=22, 0"
Replace-Text "26÷9=2, 8" "13÷2=6, 1"
Replace-Text "18÷2=9, 0" "61÷3=20, 1"
Replace-Text "86÷7=12, 2" "66÷9=7, 3"
Replace-Text "21÷2=10, 1" "79÷9=8, 7"
Replace-Text "43÷4=10, 3" "38÷2=19, 0"
Replace-Text "73÷2=36, 1" "17÷6=2, 5"
Replace-Text "98÷8=12, 2" "90÷9=10, 0"
Replace-Text "55÷2=27, 1" "47÷8=5, 7"
Replace-Text "29÷9=3, 2" "94÷4=23, 2"
Replace-Text "33÷8=4, 1" "60÷6=10, 0"
Replace-Text "21÷4=5, 1" "61÷9=6, 7"
Replace-Text "93÷9=10, 3" "66÷8=8, 2"
Replace-Text "87÷9=9, 6" "40÷7=5, 5"
Replace-Text "85÷8=10, 5" "64÷8=8, 0"
Replace-Text "60÷3=20, 0" "49÷5=9, 4"
Replace-Text "17÷2=8, 1" "31÷4=7, 3"
Replace-Text "83÷3=27, 2" "33÷7=4, 5"
Replace-Text "41÷3=13, 2" "76÷2=38, 0"
Replace-Text "24÷7=3, 3" "50÷4=12, 2"
Replace-Text "87÷7=12, 3" "15÷7=2, 1"
Replace-Text "23÷9=2, 5" "61÷6=10, 1"
